$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The Price/Volume columns store numeric-looking values (e.g. "310.98", "-0.55%")
# as plain text in the original workbook. Force text format on just the cells we
# are updating so Excel does not reinterpret them as numbers/percentages.
$textCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","D20","E20","E21","E22","D23","E23","E24","D25","E26","D39","E39","D40","E40","D41","E41","D42","E42","D44","E44","D45","E45","D46","E46","D48","E48","D50","D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '310.98'
$ws.Range('E2').Value = '-0.55%'
$ws.Range('D3').Value = '37.71'
$ws.Range('E3').Value = '-0.57%'
$ws.Range('D4').Value = '5.162'
$ws.Range('E4').Value = '1.82%'
$ws.Range('D5').Value = '0.07919'
$ws.Range('E5').Value = '1.89%'
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').Value = '1.921'
$ws.Range('E6').Value = '0.54%'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').Value = '8.287'
$ws.Range('E7').Value = '1.10%'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = '2.995'
$ws.Range('E8').Value = '0.53%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '0.9374'
$ws.Range('E9').Value = '2.01%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '0.1074'
$ws.Range('E10').Value = '-13.71%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1930'
$ws.Range('E11').Value = '1.57%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.09031'
$ws.Range('E12').Value = '0.73%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.03302'
$ws.Range('E13').Value = '-2.45%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.09618'
$ws.Range('E14').Value = '-1.06%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001378'
$ws.Range('E15').Value = '0.24%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.005718'
$ws.Range('E16').Value = '-1.85%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.597'
$ws.Range('E17').Value = '1.84%'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = '4.433'
$ws.Range('E18').Value = '1.83%'
$ws.Range('E19').Value = '0.05%'
$ws.Range('D20').Value = '6.418'
$ws.Range('E20').Value = '27.63%'
$ws.Range('E21').Value = '0.05%'
$ws.Range('E22').Value = '-2.75%'
$ws.Range('D23').Value = '0.04413'
$ws.Range('E23').Value = '0.27%'
$ws.Range('E24').Value = '1.60%'
$ws.Range('D25').Value = '0.004621'
$ws.Range('E26').Value = '0.71%'
$ws.Range('D39').Value = '0.02258'
$ws.Range('E39').Value = '5.70%'
$ws.Range('D40').Value = '0.05090'
$ws.Range('E40').Value = '2.18%'
$ws.Range('D41').Value = '0.007459'
$ws.Range('E41').Value = '-4.90%'
$ws.Range('D42').Value = '0.008888'
$ws.Range('E42').Value = '-10.18%'
$ws.Range('D44').Value = '0.002131'
$ws.Range('E44').Value = '3.37%'
$ws.Range('D45').Value = '0.009317'
$ws.Range('E45').Value = '-3.68%'
$ws.Range('D46').Value = '0.00006612'
$ws.Range('E46').Value = '1.65%'
$ws.Range('D48').Value = '0.002859'
$ws.Range('E48').Value = '-6.97%'
$ws.Range('D50').Value = '0.00002101'
$ws.Range('D51').Value = '0.0002001'
